$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Data Science Tools line: add Apache Hive, Apache Hadoop, and
#    Microsoft Azure AI to the tool list.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    " PyTorch, Sci-kit Learn, Apache Spark, MLlib, Keras, Tensorflow, LookML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " PyTorch, Sci-kit Learn, Apache Spark, Apache Hive, Apache Hadoop, MLlib, Keras, Tensorflow, LookML, Microsoft Azure AI",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the whole "Operating Systems: Mac OSX, Windows, Linux"
#    paragraph entirely (it sat right after "Backend Tools: ...").
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Operating Systems: Mac OSX, Windows, Linux`r") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) Rewrite the "Data Scientist Intern" description: the three runs
#    (sentence, space, sentence) collapse into a single run with new
#    wording and an added white shading (w:shd) on that run.
# ------------------------------------------------------------------
$oldChunkStart = "Working on Data Science project building rules-based machine learning error-detection models to carry out data quality analysis tasks on Centers for Medicare & Medicaid Services (CMS) healthcare claims data."

$find = $d.Content
$found = $find.Find.Execute("Working on Data Science*action.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Remove the matched text together with the single trailing space
    # character that immediately follows it (that trailing space run
    # will be re-created below, after the new shaded run, so ordering
    # stays correct).
    $delRng = $d.Range($find.Start, $find.End + 1)
    $delRng.Delete()

    $newText = "Working on Data Science project building rules-based machine learning error-detection models to carry out data quality analysis tasks on Centers for Medicare &amp; Medicaid Services (CMS) healthcare claims data. Building a data analysis platform using Hive, Hadoop, Apache Spark, and tools like LookML from Google's Looker. Running simulations using models to compare and test the effectiveness of different courses of action."

    $anchor1 = $d.Content
    $anchor1.Find.Execute("Data Scientist Intern", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="002B5292"><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $anchor1.InsertXML($frag1)

    $anchor2 = $d.Content
    $anchor2.Find.Execute("Data Scientist Intern", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $anchor2.InsertXML($frag2)
}

Write-Output "Done"
